$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level value changes in rows 5-25 (missing-data mask changes) ---
$ws.Range("E5").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("D19").Value = -15.5
$ws.Range("E19").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("E25").Value = -7.1

# --- Remove rows for "RM 232" (old row 26) and "SC 92" (old row 28) ---
# After deleting row 26, "SC 92" shifts up to row 27, so delete row 27 next.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Remaining per-cell value differences between the old (shifted) SC rows and the new data ---
$ws.Range("D27").ClearContents()
$ws.Range("F28").Value = 17.44
$ws.Range("E29").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("D33").Value = -14.1
